# "Seguimos ajustando autenticacion V2"
# The simulated-investment test sheet ("Datos") has its TasaEfectiva
# (effective rate) column, P, updated for both data rows from 2.89 to 4.20.
#
# The P2/P3 cells are formatted as text (NumberFormat "@") with the
# "quote prefix" cell style that was originally used when the value was
# typed as '2.89. A plain `.Value =` assignment would re-style the cell
# (dropping the quote-prefix style), so instead we stage the new text in a
# scratch cell and paste just the value into P2/P3 - this keeps the
# original cell formatting/style index intact and only changes the stored
# text, exactly like retyping the figure in Excel does.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

$scratch = $ws.Range("J6")
$scratch.Value = "4.20"
$scratch.Copy()

$ws.Range("P2").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("P3").PasteSpecial(-4163)  # xlPasteValues

$scratch.ClearContents()

$ws.Activate()
$ws.Range("P4").Select()
